# Update column O (num_edges_invaliddashed) values per the source data revision
$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("O4").Value = 0.0
$ws.Range("O6").Value = 6.0
$ws.Range("O9").Value = 0.0
$ws.Range("O11").Value = 0.0
$ws.Range("O16").Value = 1.0
$ws.Range("O17").Value = 2.0
$ws.Range("O18").Value = 3.0
$ws.Range("O19").Value = 2.0
$ws.Range("O20").Value = 1.0
$ws.Range("O22").Value = 2.0
$ws.Range("O23").Value = 0.0
$ws.Range("O24").Value = 0.0
$ws.Range("O25").Value = 0.0
$ws.Range("O27").Value = 0.0
$ws.Range("O29").Value = 2.0
$ws.Range("O32").Value = 2.0
$ws.Range("O36").Value = 1.0
$ws.Range("O37").Value = 0.0
$ws.Range("O42").Value = 2.0
$ws.Range("O43").Value = 2.0
$ws.Range("O44").Value = 3.0
$ws.Range("O46").Value = 4.0
$ws.Range("O48").Value = 0.0
$ws.Range("O50").Value = 0.0
$ws.Range("O51").Value = 4.0
$ws.Range("O52").Value = 2.0
$ws.Range("O54").Value = 0.0
$ws.Range("O55").Value = 0.0
$ws.Range("O56").Value = 0.0
$ws.Range("O58").Value = 9.0
$ws.Range("O59").Value = 1.0
$ws.Range("O60").Value = 5.0
$ws.Range("O61").Value = 0.0
$ws.Range("O62").Value = 5.0
$ws.Range("O65").Value = 0.0
$ws.Range("O67").Value = 12.0
$ws.Range("O72").Value = 2.0
$ws.Range("O75").Value = 2.0
$ws.Range("O77").Value = 4.0
$ws.Range("O81").Value = 4.0
$ws.Range("O83").Value = 1.0
$ws.Range("O87").Value = 0.0
$ws.Range("O89").Value = 2.0
$ws.Range("O90").Value = 0.0
$ws.Range("O92").Value = 6.0
